# ---------------------------------------------------------------------------
# catalogo.xlsx edit script
#
# Summary of the change being reproduced:
#   * Hoja1 renamed to "catalogo"; a new "datos" sheet is added after it.
#   * "datos" holds two helper lists: tipos (bolsos/colgantes) in column A,
#     and categoría (Chicago/Kioto/Amsterdam/Munich) in column C. Each list
#     is wrapped in its own table (Tabla2 = tipos, Tabla3 = categoría).
#   * Two new columns, "tipo" and "categoría", are inserted into the
#     products table on "catalogo" (Tabla1), between "nombre" and "precio".
#   * Two workbook-level names (tipos / categoria) point at those table
#     columns, and are used by List data-validation dropdowns added to the
#     new B2:B3 / C2:C3 ranges on "catalogo".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet and add the new "datos" sheet after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "catalogo"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "datos"

# --- 2. Make room on "catalogo" for the two new columns (tipo, categoría) --
# Inserting whole columns at B:C shifts precio/imagen1-5 from B:G to D:I and
# keeps every existing value/number-format/column-width intact.
$ws1.Range("B:C").EntireColumn.Insert()

# Grow the products table over the new range right away (before the header
# cells below are re-stamped) so it keeps tracking the header text as those
# cells are renamed.
$tbl1 = $ws1.ListObjects.Item(1)
$tbl1.Resize($ws1.Range("A1:I3"))

# --- 3. Re-stamp every header on row 1 so the table picks up the right
# column names in the right order. -------------------------------------------
$ws1.Range("C1").Value = "categoría"
$ws1.Range("B1").Value = "tipo"
$ws1.Range("D1").Value = "precio"
$ws1.Range("E1").Value = "imagen1"
$ws1.Range("F1").Value = "imagen2"
$ws1.Range("G1").Value = "imagen3"
$ws1.Range("H1").Value = "imagen4"
$ws1.Range("I1").Value = "imagen5"

# --- 4. Fill in the new column values for the two product rows -------------
$ws1.Range("C2").Value = "Chicago"
$ws1.Range("C3").Value = "Chicago"

# --- 5. Populate the "datos" sheet with the two lookup lists ----------------
$ws2.Range("A1").Value = "tipos"
$ws2.Range("A2").Value = "bolsos"
$ws2.Range("A3").Value = "colgantes"

$ws1.Range("B2").Value = "bolsos"
$ws1.Range("B3").Value = "bolsos"

$ws2.Range("C1").Value = "categoría"
$ws2.Range("C2").Value = "Chicago"
$ws2.Range("C3").Value = "Kioto"
$ws2.Range("C4").Value = "Amsterdam"
$ws2.Range("C5").Value = "Munich"

# --- 6. Create the two lookup tables (Tabla2 = tipos, Tabla3 = categoría) --
# In the source workbook these ended up numbered Tabla3 (table2.xml, id 3)
# ahead of Tabla2 (table3.xml, id 2) because table id "2" had already been
# used and freed earlier in the authoring session. Reproduce that numbering
# here with a disposable scratch table (reusing an existing shared string so
# it doesn't add a spurious sharedStrings entry) that claims id 2 and is
# removed again before the two real tables are added.
$ws2.Range("Z1").Value = "nombre"
$ws2.Range("Z2").Value = "nombre"
$tblScratch = $ws2.ListObjects.Add(1, $ws2.Range("Z1:Z2"), $null, 1)
$tblScratch.Delete()
$ws2.Range("Z1:Z2").Clear()

$tblCategoria = $ws2.ListObjects.Add(1, $ws2.Range("C1:C5"), $null, 1)
$tblCategoria.Name = "Tabla3"

$tblTipos = $ws2.ListObjects.Add(1, $ws2.Range("A1:A3"), $null, 1)
$tblTipos.Name = "Tabla2"

# --- 8. Workbook-level named ranges used by the dropdown validations -------
$wb.Names.Add("tipos", "=Tabla2[tipos]")
$wb.Names.Add("categoria", "=Tabla3[categoría]")

# --- 9. List data-validation dropdowns on the new catalogo columns ---------
$dvTipo = $ws1.Range("B2:B3")
$dvTipo.Validation.Delete()
$dvTipo.Validation.Add(3, 1, 1, "=tipos")

$dvCategoria = $ws1.Range("C2:C3")
$dvCategoria.Validation.Delete()
$dvCategoria.Validation.Add(3, 1, 1, "=categoria")

Write-Host "done"
